$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D (shifts D:K -> E:L), bringing in a new
# reporting period of financial data.
$ws.Columns.Item(4).Insert()

# The freshly inserted column D cells pick up formatting from column C by
# default; copy the formatting from column E (the old column D, now shifted
# right) onto the new column D so the new column matches the rest of the row.
# Only the rows that actually contain data are touched (36 and 78 are blank
# separator rows with no cells at all and must stay that way).
$ranges = @("E7:E35", "E38:E77", "E80:E102")
foreach ($rng in $ranges) {
    $src = $ws.Range($rng)
    $col = $src.Column
    $dst = $ws.Range($src.Address()).Offset(0, -1)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Populate the new column D with the new period's values.
$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(8,4).Value = 1804900
$ws.Cells.Item(9,4).Value = 709500
$ws.Cells.Item(10,4).Value = 1095400
$ws.Cells.Item(12,4).Value = "NA"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(14,4).Value = 20300
$ws.Cells.Item(15,4).Value = 90900
$ws.Cells.Item(17,4).Value = 1735600
$ws.Cells.Item(18,4).Value = 69300
$ws.Cells.Item(20,4).Value = -500
$ws.Cells.Item(21,4).Value = 159700
$ws.Cells.Item(22,4).Value = 34900
$ws.Cells.Item(23,4).Value = 34000
$ws.Cells.Item(24,4).Value = 7900
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(26,4).Value = 26100
$ws.Cells.Item(27,4).Value = 24900
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(32,4).Value = 500
$ws.Cells.Item(33,4).Value = 24900
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(35,4).Value = 24900
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(41,4).Value = 9900
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(43,4).Value = 210700
$ws.Cells.Item(44,4).Value = 7300
$ws.Cells.Item(45,4).Value = 57600
$ws.Cells.Item(46,4).Value = 285500
$ws.Cells.Item(47,4).Value = "NA"
$ws.Cells.Item(48,4).Value = 518700
$ws.Cells.Item(49,4).Value = 86600
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(52,4).Value = 19600
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(54,4).Value = 910500
$ws.Cells.Item(57,4).Value = 63800
$ws.Cells.Item(58,4).Value = 113100
$ws.Cells.Item(59,4).Value = 87400
$ws.Cells.Item(60,4).Value = 264300
$ws.Cells.Item(61,4).Value = 311500
$ws.Cells.Item(62,4).Value = 88000
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(66,4).Value = 675600
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(72,4).Value = -17300
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(76,4).Value = 234900
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(81,4).Value = 24900
$ws.Cells.Item(83,4).Value = 90800
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(89,4).Value = 100600
$ws.Cells.Item(91,4).Value = -223900
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(94,4).Value = -166100
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(100,4).Value = 66200
$ws.Cells.Item(101,4).Value = 0
$ws.Cells.Item(102,4).Value = 700
